$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100.7
$ws.Range("I5").Value = 84.875
$ws.Range("K5").Value = 84.875
$ws.Range("M5").Value = 30.125
$ws.Range("H9").Value = 7678.5
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 21035.5
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 21035.5
$ws.Range("M9").Value = -831
$ws.Range("N9").Value = -21373.5
$ws.Range("H39").Value = 691.7778
$ws.Range("I39").Value = 79.5
$ws.Range("K39").Value = 238.5
$ws.Range("M39").Value = 57.5
$ws.Range("H64").Value = 8658
$ws.Range("I64").Value = 7350
$ws.Range("J64").Value = 9966
$ws.Range("K64").Value = 7350
$ws.Range("L64").Value = 9966
$ws.Range("M64").Value = -7102
$ws.Range("N64").Value = -10462
$ws.Range("H67").Value = 8658
$ws.Range("I67").Value = 7350
$ws.Range("J67").Value = 9966
$ws.Range("K67").Value = 7350
$ws.Range("L67").Value = 9966
$ws.Range("M67").Value = -6492
$ws.Range("N67").Value = -11682
$ws.Range("H97").Value = 3613.125
$ws.Range("J97").Value = 3613.125
$ws.Range("L97").Value = 10839.375
$ws.Range("N97").Value = -11831.375
$ws.Range("H107").Value = 1071.1111
$ws.Range("I107").Value = 1243.6
$ws.Range("J107").Value = 208.66667
$ws.Range("K107").Value = 1243.6
$ws.Range("L107").Value = 208.66667
$ws.Range("M107").Value = 676.4000000000001
$ws.Range("N107").Value = -4048.66667
$ws.Range("H112").Value = 77024.86
$ws.Range("J112").Value = 89196.164
$ws.Range("L112").Value = 267588.492
$ws.Range("N112").Value = -269804.492
$ws.Range("H132").Value = 1296.1428
$ws.Range("I132").Value = 1224.6666
$ws.Range("K132").Value = 3673.9998
$ws.Range("M132").Value = -1143.9998
$ws.Range("H137").Value = 2389220.8
$ws.Range("I137").Value = 7236.2334
$ws.Range("J137").Value = 8344182
$ws.Range("K137").Value = 21708.7002
$ws.Range("L137").Value = 25032546
$ws.Range("M137").Value = -19158.7002
$ws.Range("N137").Value = -25037646

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 220756.44
$ws.Range("I32").Value = 278828.8
$ws.Range("J32").Value = 11695.9
$ws.Range("K32").Value = 278828.8
$ws.Range("L32").Value = 11695.9
$ws.Range("M32").Value = -278541.8
$ws.Range("N32").Value = -12269.9
$ws.Range("H39").Value = 12502500
$ws.Range("I39").Value = 12502500
$ws.Range("K39").Value = 12502500
$ws.Range("M39").Value = -12501980
$ws.Range("H61").Value = 1464791.9
$ws.Range("I61").Value = 50699.637
$ws.Range("K61").Value = 50699.637
$ws.Range("M61").Value = -50487.637
$ws.Range("H74").Value = 485277.62
$ws.Range("I74").Value = 1614.975
$ws.Range("K74").Value = 1614.975
$ws.Range("M74").Value = -740.9749999999999
$ws.Range("H77").Value = 485277.62
$ws.Range("I77").Value = 1614.975
$ws.Range("K77").Value = 8074.875
$ws.Range("M77").Value = -3706.875
$ws.Range("H132").Value = 2886.9167
$ws.Range("I132").Value = 2644.3635
$ws.Range("J132").Value = 5555
$ws.Range("K132").Value = 7933.0905
$ws.Range("L132").Value = 16665
$ws.Range("M132").Value = -5403.0905
$ws.Range("N132").Value = -21725
$ws.Range("H136").Value = 1464791.9
$ws.Range("I136").Value = 50699.637
$ws.Range("K136").Value = 152098.911
$ws.Range("M136").Value = -149548.911

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 52497.5
$ws.Range("I93").Value = 49995
$ws.Range("J93").Value = 55000
$ws.Range("K93").Value = 49995
$ws.Range("L93").Value = 55000
$ws.Range("M93").Value = -48123
$ws.Range("N93").Value = -58744
$ws.Range("H99").Value = 6377.2
$ws.Range("I99").Value = 6896.8076
$ws.Range("K99").Value = 6896.8076
$ws.Range("M99").Value = -5398.8076
$ws.Range("H107").Value = 19646.572
$ws.Range("J107").Value = 8657.833000000001
$ws.Range("L107").Value = 8657.833000000001
$ws.Range("N107").Value = -12497.833
$ws.Range("H134").Value = 29034102
$ws.Range("I134").Value = 1922.3462
$ws.Range("K134").Value = 5767.0386
$ws.Range("M134").Value = -3232.0386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5198.8647
$ws.Range("I31").Value = 4271.793
$ws.Range("K31").Value = 4271.793
$ws.Range("M31").Value = -3976.793
$ws.Range("H34").Value = 5198.8647
$ws.Range("I34").Value = 4271.793
$ws.Range("K34").Value = 4271.793
$ws.Range("M34").Value = -4069.793
$ws.Range("H45").Value = 49850
$ws.Range("J45").Value = 49850
$ws.Range("L45").Value = 49850
$ws.Range("N45").Value = -51036
$ws.Range("H132").Value = 2736.95
$ws.Range("I132").Value = 1835.091
$ws.Range("K132").Value = 5505.272999999999
$ws.Range("M132").Value = -2975.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9969.823
$ws.Range("I3").Value = 8706.714
$ws.Range("J3").Value = 15864.333
$ws.Range("K3").Value = 26120.142
$ws.Range("L3").Value = 47592.999
$ws.Range("M3").Value = -26008.142
$ws.Range("N3").Value = -47816.999
$ws.Range("H5").Value = 1149.8
$ws.Range("I5").Value = 687.25
$ws.Range("K5").Value = 2061.75
$ws.Range("M5").Value = -1949.75
$ws.Range("H107").Value = 668.125
$ws.Range("I107").Value = 589.0526
$ws.Range("J107").Value = 739.6667
$ws.Range("K107").Value = 1767.1578
$ws.Range("L107").Value = 2219.0001
$ws.Range("M107").Value = 152.8422
$ws.Range("N107").Value = -6059.0001
$ws.Range("H126").Value = 12666.167
$ws.Range("J126").Value = 14199.4
$ws.Range("L126").Value = 42598.2
$ws.Range("N126").Value = -52478.2
$ws.Range("H135").Value = 1149.8
$ws.Range("I135").Value = 687.25
$ws.Range("K135").Value = 6185.25
$ws.Range("M135").Value = -3650.25
$ws.Range("H137").Value = 3936
$ws.Range("J137").Value = 9999.5
$ws.Range("L137").Value = 29998.5
$ws.Range("N137").Value = -40198.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4683.6924
$ws.Range("J126").Value = 5286
$ws.Range("L126").Value = 15858
$ws.Range("N126").Value = -20798
$ws.Range("H132").Value = 5701606
$ws.Range("I132").Value = 3666918.8
$ws.Range("J132").Value = 12823012
$ws.Range("K132").Value = 11000756.4
$ws.Range("L132").Value = 38469036
$ws.Range("M132").Value = -10998226.4
$ws.Range("N132").Value = -38474096

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7840.476
$ws.Range("J7").Value = 12999.111
$ws.Range("L7").Value = 12999.111
$ws.Range("N7").Value = -13223.111
$ws.Range("H122").Value = 3003.2188
$ws.Range("I122").Value = 2754.4138
$ws.Range("K122").Value = 8263.241399999999
$ws.Range("M122").Value = -5813.241399999999
$ws.Range("H124").Value = 97429
$ws.Range("J124").Value = 97429
$ws.Range("L124").Value = 97429
$ws.Range("N124").Value = -107249
$ws.Range("H125").Value = 91483
$ws.Range("J125").Value = 91483
$ws.Range("L125").Value = 91483
$ws.Range("N125").Value = -101323
$ws.Range("H126").Value = 7840.476
$ws.Range("J126").Value = 12999.111
$ws.Range("L126").Value = 38997.333
$ws.Range("N126").Value = -43937.333
$ws.Range("H127").Value = 91974.75
$ws.Range("J127").Value = 91974.75
$ws.Range("L127").Value = 91974.75
$ws.Range("N127").Value = -101894.75
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 35000
$ws.Range("J129").Value = 35000
$ws.Range("L129").Value = 35000
$ws.Range("M129").Value = -45000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 20000
$ws.Range("K49").Value = 20000
$ws.Range("M49").Value = -19770
$ws.Range("H136").Value = 2862.1724
$ws.Range("I136").Value = 2393.4348
$ws.Range("K136").Value = 7180.3044
$ws.Range("M136").Value = -4630.3044
